# Commit: "Changed worksheet name for mif_template and micsss_template (#82)"
#
# This mIF template workbook had its first worksheet renamed from
# "mIF_template" to "mIF", a couple of data-entry values corrected
# (slide scanner model, autostainer model, protocol name typo), and the
# last-used selection on each sheet updated to reflect where the editor
# left the cursor.

$wb = $excel.ActiveWorkbook

$wsMif  = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# 1) Rename the main worksheet tab: "mIF_template" -> "mIF"
$wsMif.Name = "mIF"

# 2) Fix up a few values in the template's sample/header rows.
#    C5  = SLIDE SCANNER MODEL   : "Vectra 2.0"    -> "Hamamatsu"
#    C7  = AUTOSTAINER MODEL     : "Leica Bon RX"   -> "Leica Bond RX" (typo fix)
#    C12 = PROTOCOL NAME         : "T-Cell HNSC"    -> "T-Cell HSNC"   (typo fix)
$wsMif.Range("C5").Value = "Hamamatsu"
$wsMif.Range("C7").Value = "Leica Bond RX"
$wsMif.Range("C12").Value = "T-Cell HSNC"

# 3) Leave the cursor/selection where the editor ended up on each sheet.
[void]$wsMif.Activate()
[void]$wsMif.Range("B16").Select()

[void]$wsData.Activate()
[void]$wsData.Range("I5").Select()

# End back on the main worksheet, which is the tab that should be
# active/selected when the workbook is reopened.
[void]$wsMif.Activate()
